$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1740
$ws.Range("F4").Value = 129
$ws.Range("F5").Value = 356
$ws.Range("F6").Value = 756
$ws.Range("F7").Value = 198
$ws.Range("C8").Value = "上海·幻想物语新春动漫嘉年华"
$ws.Range("D8").Value = "长寿路309号 旭辉企业大厦"
$ws.Range("E8").Value = "2024.02.16 10:00-02.17 17:00"
$ws.Range("F8").Value = 1062
$ws.Range("G8").Value = 68
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=81682"
$ws.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202402/bsBJZToU1707285951479.jpeg"
$ws.Range("C9").Value = "上海·次元裂缝-X 新年anikura派对"
$ws.Range("D9").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("E9").Value = "2024.02.16 14:00-02.16 19:00"
$ws.Range("F9").Value = 263
$ws.Range("G9").Value = 60
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=81314"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202401/OrhHWKdR1706062360956.jpeg"
$ws.Range("F11").Value = 347
$ws.Range("F12").Value = 614
$ws.Range("F17").Value = 153
$ws.Range("F18").Value = 818
$ws.Range("F25").Value = 16
$ws.Range("F27").Value = 583
$ws.Range("F28").Value = 956
$ws.Range("F29").Value = 1
$ws.Range("F31").Value = 217
$ws.Range("F32").Value = 1019
$ws.Range("F34").Value = 42

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1035
$ws.Range("F5").Value = 1035
$ws.Range("F10").Value = 320
$ws.Range("F14").Value = 578
$ws.Range("F15").Value = 90
$ws.Range("F17").Value = 967
$ws.Range("F22").Value = 33
$ws.Range("F24").Value = 289
$ws.Range("F25").Value = 259
$ws.Range("F26").Value = 3601
$ws.Range("G26").Value = "已售罄"
$ws.Range("F31").Value = 18
$ws.Range("F33").Value = 101
$ws.Range("F34").Value = 26

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1767
$ws.Range("F4").Value = 42
$ws.Range("F5").Value = 2395
$ws.Range("F6").Value = 984
$ws.Range("F9").Value = 1238
$ws.Range("F10").Value = 324

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1767
$ws.Range("F3").Value = 42
$ws.Range("F4").Value = 2395
$ws.Range("F6").Value = 1740
$ws.Range("F7").Value = 984
$ws.Range("F8").Value = 1238
$ws.Range("F9").Value = 324
$ws.Range("F11").Value = 129
$ws.Range("F12").Value = 356
$ws.Range("F13").Value = 756
$ws.Range("F14").Value = 198
$ws.Range("C16").Value = "上海·幻想物语新春动漫嘉年华"
$ws.Range("D16").Value = "长寿路309号 旭辉企业大厦"
$ws.Range("E16").Value = "2024.02.16 10:00-02.17 17:00"
$ws.Range("F16").Value = 1062
$ws.Range("G16").Value = 68
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=81682"
$ws.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202402/bsBJZToU1707285951479.jpeg"
$ws.Range("C17").Value = "上海·次元裂缝-X 新年anikura派对"
$ws.Range("D17").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("E17").Value = "2024.02.16 14:00-02.16 19:00"
$ws.Range("F17").Value = 263
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=81314"
$ws.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202401/OrhHWKdR1706062360956.jpeg"
$ws.Range("F18").Value = 347
$ws.Range("F19").Value = 614
$ws.Range("F20").Value = 1035
$ws.Range("F23").Value = 153
$ws.Range("F24").Value = 818
$ws.Range("F29").Value = 320
$ws.Range("F33").Value = 583
$ws.Range("F34").Value = 956
$ws.Range("F35").Value = 578
$ws.Range("F36").Value = 578
$ws.Range("F37").Value = 90
$ws.Range("F39").Value = 217
$ws.Range("F43").Value = 33
$ws.Range("F44").Value = 289
$ws.Range("F45").Value = 289
$ws.Range("F46").Value = 259
$ws.Range("F47").Value = 1019
$ws.Range("F50").Value = 42
